# Apply the "Dep Ed Closures" update:
#  1. Update the "as at" date/time line in A7.
#  2. Insert a new early-childhood-service closure entry before the
#     existing "VERMONT CHILDREN'S CENTRE" row (North-Eastern region list).
#  3. Insert a new early-childhood-service closure entry before the
#     existing "Sunshine Leisure Centre" row (South-Western region list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the date/time summary text in A7.
$oldA7 = "On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school" + [char]160 + "and early childhood service," + [char]160 + "TAFE closures and relocations for Monday 31 Augst," + [char]160 + "(as at 11:10am, 31August)South-Eastern Victoria RegionEarly childhood services"
$newA7 = "On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school" + [char]160 + "and early childhood service," + [char]160 + "TAFE closures and relocations for Tuesday 1 September," + [char]160 + "(as at 10:30pm, 31August)South-Eastern Victoria RegionEarly childhood services"
$ws.Range("A7").Value = $newA7

# 2. Insert new row 111 (pushes the VERMONT CHILDREN'S CENTRE row, and everything
#    after it, down by one row) and populate it.
$ws.Rows.Item(111).Insert()
$ws.Range("A111").Value = "li: Toolamba Outside School Hours Care TOOLAMBA"

# 3. Insert new row 204 (pushes the Sunshine Leisure Centre row, and everything
#    after it, down by one more row) and populate it.
$ws.Rows.Item(204).Insert()
$ws.Range("A204").Value = "li: South Pacific Health Club Williamstown NEWPORT"
